$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the promotional text with new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$oldText = $ws1.Range("A1").Value2
$newText = $oldText `
    -replace [regex]::Escape("1000 Bs = 3.34 = 13006.69 pesos"), "1000 Bs = 3.31 = 12868.01 pesos" `
    -replace [regex]::Escape("13006.69 pesos = 3.32 = 964.23 Bs"), "12868.01 pesos = 3.3 = 964.11 Bs"
$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the N10/O10/N12/O12 rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 301.99
$ws2.Range("O10").Value = 3886.01
$ws2.Range("N12").Value = 3900
$ws2.Range("O12").Value = 292.2
